# Re-colour the deck's theme from the custom "Integral" (Red Violet) palette
# to the stock Office Theme palette, and point the three data tables at the
# corresponding built-in Office table style (instead of the custom one that
# shipped with the old theme).

$p = $ppt.ActivePresentation

# --- 1. Theme colours -------------------------------------------------
# Slide.ThemeColorScheme exposes the 12 theme colour slots in the fixed
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink. Re-pointing them
# at the Office Theme RGB values updates the presentation's single shared
# theme part used by every slide/layout/master.
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1
$tcs.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink -> 954F72

# --- 2. Table styles ----------------------------------------------------
# The three tables (on the slides that hold them) move off the custom
# "Table_0" style that belonged to the old theme and onto the standard
# Office table style.
$newTableStyle = "{AD235E05-E29A-404F-AB7B-ABAA8784064B}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}
